# Update cryptocurrency price/volume data per upstream refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '34.419.22'
$ws.Range('E2').Value = "'" + '  +0.88%  '
$ws.Range('D3').Value = "'" + '1.794.35'
$ws.Range('E3').Value = "'" + '  +0.36%  '
$ws.Range('E4').Value = "'" + '  +0.11%  '
$ws.Range('D5').Value = "'" + '226.79'
$ws.Range('E6').Value = "'" + '  +1.50%  '
$ws.Range('E7').Value = "'" + '  +0.11%  '
$ws.Range('D8').Value = "'" + '32.36'
$ws.Range('E8').Value = "'" + '  +1.26%  '
$ws.Range('D9').Value = "'" + '0.295'
$ws.Range('E9').Value = "'" + '  +1.17%  '
$ws.Range('D10').Value = "'" + '0.0693'
$ws.Range('E10').Value = "'" + '  +0.40%  '
$ws.Range('D11').Value = "'" + '0.0950'
$ws.Range('E11').Value = "'" + '  +0.54%  '
$ws.Range('D12').Value = "'" + '2.055.79'
$ws.Range('E12').Value = "'" + '  +0.49%  '
$ws.Range('D13').Value = "'" + '11.03'
$ws.Range('E13').Value = "'" + '  -1.40%  '
$ws.Range('D14').Value = "'" + '1.789.40'
$ws.Range('E14').Value = "'" + '  +0.06%  '
$ws.Range('D15').Value = "'" + '0.631'
$ws.Range('E15').Value = "'" + '  +1.80%  '
$ws.Range('D16').Value = "'" + '34.387.62'
$ws.Range('E16').Value = "'" + '  +0.94%  '
$ws.Range('D17').Value = "'" + '4.22'
$ws.Range('E17').Value = "'" + '  +1.07%  '
$ws.Range('D18').Value = "'" + '68.30'
$ws.Range('E18').Value = "'" + '  +0.36%  '
$ws.Range('D19').Value = "'" + '0.0₃0802'
$ws.Range('E19').Value = "'" + '  +3.04%  '
$ws.Range('D20').Value = "'" + '246.38'
$ws.Range('E20').Value = "'" + '  +0.46%  '
$ws.Range('E22').Value = "'" + '  -0.01%  '
$ws.Range('D23').Value = "'" + '4.14'
$ws.Range('E23').Value = "'" + '  +1.04%  '
$ws.Range('D25').Value = "'" + '162.88'
$ws.Range('E25').Value = "'" + '  +0.94%  '
$ws.Range('D26').Value = "'" + '7.23'
$ws.Range('E26').Value = "'" + '  +1.12%  '
$ws.Range('D27').Value = "'" + '16.39'
$ws.Range('E27').Value = "'" + '  +0.50%  '
$ws.Range('E28').Value = "'" + '  +2.24%  '
$ws.Range('E29').Value = "'" + '  +0.14%  '
$ws.Range('E30').Value = "'" + '  +0.70%  '
$ws.Range('D31').Value = "'" + '0.0521'
$ws.Range('E31').Value = "'" + '  +0.94%  '
$ws.Range('B32').Value = "'" + 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = "'" + 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').Value = "'" + '3.89'
$ws.Range('E32').Value = "'" + '  +8.17%  '
$ws.Range('B33').Value = "'" + 'Filecoin'
$ws.Range('C33').Value = "'" + 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = "'" + '3.77'
$ws.Range('E33').Value = "'" + '  +3.18%  '
$ws.Range('E34').Value = "'" + '  +1.23%  '
$ws.Range('D35').Value = "'" + '1.443.88'
$ws.Range('E35').Value = "'" + '  -0.92%  '
$ws.Range('D36').Value = "'" + '2.63'
$ws.Range('E36').Value = "'" + '  +8.77%  '
$ws.Range('D37').Value = "'" + '0.667'
$ws.Range('E37').Value = "'" + '  +3.27%  '
$ws.Range('E38').Value = "'" + '  +1.84%  '
$ws.Range('E39').Value = "'" + '  -0.77%  '
$ws.Range('D40').Value = "'" + '83.91'
$ws.Range('E40').Value = "'" + '  +4.46%  '
$ws.Range('E41').Value = "'" + '  +1.45%  '
$ws.Range('B42').Value = "'" + 'ARBITRUM'
$ws.Range('C42').Value = "'" + 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D42').Value = "'" + '0.934'
$ws.Range('E42').Value = "'" + '  +1.79%  '
$ws.Range('B43').Value = "'" + 'MXToken'
$ws.Range('C43').Value = "'" + 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').Value = "'" + '2.75'
$ws.Range('E43').Value = "'" + '  +2.70%  '
$ws.Range('D44').Value = "'" + '13.80'
$ws.Range('E44').Value = "'" + '  +2.14%  '
$ws.Range('E45').Value = "'" + '  +3.96%  '
$ws.Range('E46').Value = "'" + '  +0.62%  '
$ws.Range('D47').Value = "'" + '1.08'
$ws.Range('E47').Value = "'" + '  +0.24%  '
$ws.Range('D48').Value = "'" + '1.952.07'
$ws.Range('E48').Value = "'" + '  +0.27%  '
$ws.Range('D49').Value = "'" + '105.69'
$ws.Range('E49').Value = "'" + '  -1.02%  '
$ws.Range('B50').Value = "'" + 'BabyDogeCoin'
$ws.Range('C50').Value = "'" + 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = "'" + '0.0₆0131'
$ws.Range('E50').Value = "'" + '  -2.72%  '
$ws.Range('B51').Value = "'" + 'PaxDollar'
$ws.Range('C51').Value = "'" + 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D51').Value = "'" + '1.00'
$ws.Range('E51').Value = "'" + '  +0.11%  '
